# Generate Report for Handoff
#
# The source file's GUID changed from 37d393ee-43ce-439f-8ab2-38282a781ca7
# to e2bfda48-2bf0-4296-8247-59b7b7126d9d, and a new handoff round produced
# new handoff-file names + handoff timestamps for both locales. The
# underlying hyperlink targets (the relationship URLs) are left exactly as
# they were -- only the cell text / hyperlink display text is refreshed to
# reflect the new file names and timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "37d393ee-43ce-439f-8ab2-38282a781ca7"
$newGuid = "e2bfda48-2bf0-4296-8247-59b7b7126d9d"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$oldZhXlf = "$oldGuid.1f0e9e31cf49d3412ace53f43ee3c98a140be73d.zh-cn.xlf"
$newZhXlf = "$newGuid.6402398c32532f0b25ef07c097b6cf1ca93d684c.zh-cn.xlf"

$oldDeXlf = "$oldGuid.1f0e9e31cf49d3412ace53f43ee3c98a140be73d.de-de.xlf"
$newDeXlf = "$newGuid.6402398c32532f0b25ef07c097b6cf1ca93d684c.de-de.xlf"

$oldZhTime = "2016-03-09 01:27:21"
$newZhTime = "2016-03-09 01:28:11"

$oldDeTime = "2016-03-09 01:27:31"
$newDeTime = "2016-03-09 01:28:21"

# Hyperlink targets (relationship URLs) -- unchanged by this edit, reused
# verbatim when the hyperlinks are re-created below.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/99c46e84003a1186abeae2ddc6084f3e3b08a06a/e2e/$oldMdName"
$configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/99c46e84003a1186abeae2ddc6084f3e3b08a06a/.localization-config"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a29d120e6782f9a653d8efac6e56e799709f1a02/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3c2571b70b0c6bc520fb418b9737cf809a0c58f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# --- Sheet "Overview": only A2 (the .md hyperlink) mentions the GUID -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddress, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configAddress, "", "", ".localization-config")

# --- Sheet "zh-cn": A2 (.md), C2 (.xlf handoff file), D2 (handoff time) --
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddress, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfAddress, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configAddress, "", "", ".localization-config")
$wsZh.Range("D2").Value = $newZhTime

# --- Sheet "de-de": A2 (.md), C2 (.xlf handoff file), D2 (handoff time) --
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddress, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfAddress, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configAddress, "", "", ".localization-config")
$wsDe.Range("D2").Value = $newDeTime
